$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each (cell, new text value) pair below mirrors the diff exactly.
# We force the Text number format before assignment so Excel keeps the
# exact literal string (matching the original inlineStr cells) instead of
# re-interpreting it as a number/percentage, then reset the style so no
# extra formatting is left on the cell (matches the unstyled source cells).
$cellValues = [ordered]@{
    'D2' = '258.47'
    'E2' = '0.84%'
    'D3' = '26.84'
    'E3' = '-1.61%'
    'D4' = '4.643'
    'E4' = '0.67%'
    'D5' = '0.05953'
    'E5' = '1.07%'
    'D6' = '6.628'
    'E6' = '-0.25%'
    'D7' = '0.8563'
    'E7' = '-0.93%'
    'D8' = '0.9272'
    'E8' = '-0.22%'
    'D9' = '0.1386'
    'E9' = '-1.60%'
    'D10' = '0.04373'
    'E10' = '16.54%'
    'D11' = '0.07007'
    'E11' = '-1.17%'
    'D12' = '0.02969'
    'E12' = '-8.01%'
    'D13' = '0.09112'
    'E13' = '-1.02%'
    'D14' = '0.001533'
    'E14' = '-0.17%'
    'D15' = '0.0006038'
    'E15' = '-0.19%'
    'D16' = '0.006114'
    'E16' = '0.08%'
    'D17' = '3.452'
    'E17' = '-1.77%'
    'D18' = '3.127'
    'E18' = '-1.98%'
    'E19' = '-2.19%'
    'E20' = '-0.01%'
    'D21' = '0.1297'
    'E21' = '1.69%'
    'D22' = '3.827'
    'E22' = '-0.90%'
    'D23' = '0.04219'
    'E23' = '-0.18%'
    'E24' = '-0.55%'
    'E25' = '4.67%'
    'D26' = '0.0001199'
    'E26' = '-0.12%'
    'E27' = '-11.63%'
    'D40' = '0.03816'
    'E40' = '-0.20%'
    'D41' = '0.1108'
    'E41' = '0.95%'
    'D42' = '0.003779'
    'E42' = '-39.64%'
    'E43' = '27.73%'
    'D44' = '0.01491'
    'E44' = '31.26%'
    'D45' = '0.00005132'
    'E45' = '-6.40%'
    'D46' = '0.00000000750'
    'D47' = '0.04997'
    'E47' = '-16.98%'
    'D48' = '0.2199'
    'E48' = '9,544.46%'
    'D49' = '0.00002099'
    'D50' = '0.0001999'
}

foreach ($cell in $cellValues.Keys) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $cellValues[$cell]
    $rng.Style = "Normal"
}
